$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-42): refreshed weekly price/volume/date data ---
$ws.Range("D2").Value = 44350
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 24000
$ws.Range("M2").Value = 24375
$ws.Range("P2").Value = 2438

$ws.Range("D3").Value = 44356
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = 24000
$ws.Range("P3").Value = 2400

$ws.Range("D4").Value = 44349
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 24000
$ws.Range("P4").Value = 2400

$ws.Range("D5").Value = 44389
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("P5").Value = 2500

$ws.Range("D6").Value = 44382
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("P6").Value = 2500

$ws.Range("D7").Value = 44390
$ws.Range("J7").Value = 15

$ws.Range("D9").Value = 44386
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 25000
$ws.Range("P9").Value = 2500

$ws.Range("D10").Value = 44354
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 24000
$ws.Range("L10").Value = 24000
$ws.Range("M10").Value = 24000
$ws.Range("P10").Value = 2400

$ws.Range("D11").Value = 44413
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 25000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 25000
$ws.Range("P11").Value = 2500

$ws.Range("D12").Value = 44405
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 25000
$ws.Range("P12").Value = 2500

$ws.Range("D13").Value = 44431
$ws.Range("J13").Value = 65

$ws.Range("D14").Value = 44355
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = 23000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23400
$ws.Range("P14").Value = 2340

$ws.Range("D15").Value = 44372
$ws.Range("J15").Value = 20

$ws.Range("D16").Value = 44410
$ws.Range("J16").Value = 50

$ws.Range("D17").Value = 44396
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 25000
$ws.Range("P17").Value = 2500

$ws.Range("D18").Value = 44406
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 24000
$ws.Range("M18").Value = 24542
$ws.Range("P18").Value = 2454

$ws.Range("D19").Value = 44385
$ws.Range("J19").Value = 80

$ws.Range("D20").Value = 44371
$ws.Range("J20").Value = 50

$ws.Range("D21").Value = 44419
$ws.Range("J21").Value = 25

$ws.Range("D22").Value = 44392
$ws.Range("J22").Value = 25

$ws.Range("D23").Value = 44384
$ws.Range("J23").Value = 40

$ws.Range("D24").Value = 44365
$ws.Range("J24").Value = 85
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 22000
$ws.Range("P24").Value = 2200

$ws.Range("D25").Value = 44433
$ws.Range("J25").Value = 25
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 25000
$ws.Range("P25").Value = 2500

$ws.Range("D26").Value = 44397
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = 27000
$ws.Range("L26").Value = 27000
$ws.Range("M26").Value = 27000
$ws.Range("P26").Value = 2700

$ws.Range("D27").Value = 44426
$ws.Range("J27").Value = 30

$ws.Range("D28").Value = 44421
$ws.Range("J28").Value = 55

$ws.Range("D29").Value = 44434
$ws.Range("J29").Value = 55
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 25000
$ws.Range("P29").Value = 2500

$ws.Range("D30").Value = 44420
$ws.Range("J30").Value = 55

$ws.Range("D31").Value = 44348
$ws.Range("J31").Value = 3

$ws.Range("D32").Value = 44427
$ws.Range("J32").Value = 40

$ws.Range("D33").Value = 44417
$ws.Range("J33").Value = 15
$ws.Range("K33").Value = 25000
$ws.Range("M33").Value = 25000
$ws.Range("P33").Value = 2500

$ws.Range("D34").Value = 44441
$ws.Range("J34").Value = 70

$ws.Range("D35").Value = 44432
$ws.Range("K35").Value = 27000
$ws.Range("L35").Value = 27000
$ws.Range("M35").Value = 27000
$ws.Range("P35").Value = 2700

$ws.Range("D36").Value = 44379
$ws.Range("J36").Value = 35
$ws.Range("K36").Value = 22000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 22000
$ws.Range("P36").Value = 2200

$ws.Range("D37").Value = 44446
$ws.Range("J37").Value = 40
$ws.Range("K37").Value = 27000
$ws.Range("L37").Value = 27000
$ws.Range("M37").Value = 27000
$ws.Range("P37").Value = 2700

$ws.Range("D38").Value = 44411
$ws.Range("J38").Value = 40

$ws.Range("D39").Value = 44438

$ws.Range("D40").Value = 44428
$ws.Range("J40").Value = 30
$ws.Range("K40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 25000
$ws.Range("P40").Value = 2500

$ws.Range("D41").Value = 44435
$ws.Range("J41").Value = 185
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 27000
$ws.Range("M41").Value = 25162
$ws.Range("P41").Value = 2516

$ws.Range("D42").Value = 44376
$ws.Range("J42").Value = 45
$ws.Range("K42").Value = 23000
$ws.Range("L42").Value = 23000
$ws.Range("M42").Value = 23000
$ws.Range("P42").Value = 2300

# --- Append new rows 43-44 (additional weekly observations) ---
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44412
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100112035
$ws.Range("G43").Value = "Bruselas (repollito)"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 50
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = 25000
$ws.Range("N43").Value = "`$/malla 10 kilos"
$ws.Range("O43").Value = "Provincia de Quillota"
$ws.Range("P43").Value = 2500
$ws.Range("Q43").Value = 10
$ws.Range("R43").Value = "Hortaliza"

$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44400
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112035
$ws.Range("G44").Value = "Bruselas (repollito)"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 12
$ws.Range("K44").Value = 24000
$ws.Range("L44").Value = 24000
$ws.Range("M44").Value = 24000
$ws.Range("N44").Value = "`$/malla 10 kilos"
$ws.Range("O44").Value = "Provincia de Quillota"
$ws.Range("P44").Value = 2400
$ws.Range("Q44").Value = 10
$ws.Range("R44").Value = "Hortaliza"

Write-Host "Applied weekly update to Bruselas (repollito) sheet"